$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A51").Value = "2025-04-29 05:55:17"
$ws.Range("B51").Value = 163
